$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-5, columns B:E with the new cluster analysis values.
$ws.Range("B2").Value = 86
$ws.Range("C2").Value = -4.958607314841775
$ws.Range("D2").Value = -0.9866193802049792
$ws.Range("E2").Value = -2.640283785166799

$ws.Range("B3").Value = 151
$ws.Range("C3").Value = -5.527243550682788
$ws.Range("D3").Value = -1.220981028085129
$ws.Range("E3").Value = -2.736433256119832

$ws.Range("B4").Value = 161
$ws.Range("C4").Value = -5.376750709602099
$ws.Range("D4").Value = -0.08092190762392611
$ws.Range("E4").Value = -2.323090469396636

$ws.Range("B5").Value = 48
$ws.Range("C5").Value = -5.283996656365201
$ws.Range("D5").Value = -1.356547323513813
$ws.Range("E5").Value = -3.261739112424002
